$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "fylker_oppdatert" (C) and "fylker_ltmv" (D) columns for the
# existing rows to reflect the new fylke (county) division that took effect
# from 2020 ("nyeste fylkesinndeling").

$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 11

$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 11

$ws.Range("D4").Value = 10

$ws.Range("C5").Value = 34
$ws.Range("D5").Value = 9

$ws.Range("C6").Value = 34
$ws.Range("D6").Value = 9

$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 11

$ws.Range("C8").Value = 38
$ws.Range("D8").Value = 8

$ws.Range("C9").Value = 38
$ws.Range("D9").Value = 8

$ws.Range("C10").Value = 42
$ws.Range("D10").Value = 7

$ws.Range("C11").Value = 42
$ws.Range("D11").Value = 7

$ws.Range("D12").Value = 6

$ws.Range("C13").Value = 46
$ws.Range("D13").Value = 5

$ws.Range("C14").Value = 46
$ws.Range("D14").Value = 5

$ws.Range("D15").Value = 4

$ws.Range("D16").Value = 3

$ws.Range("D17").Value = 3

$ws.Range("D18").Value = 2

$ws.Range("C19").Value = 54
$ws.Range("D19").Value = 1

$ws.Range("C20").Value = 54

$ws.Range("D22").Value = 3

# --- Append the six new counties (fylker) introduced by the reform as new
# rows 23-28, with their own code in fylker_oppdatert/fylker_ltmv and the
# associated RHF (regional health authority) info.

$newRows = @(
    @{ Row = 23; A = 30; B = "Viken";                C = 30; D = 11; E = 1; F = "Helse Sør-Øst"; G = 111919 },
    @{ Row = 24; A = 34; B = "Innlandet";             C = 34; D = 9;  E = 1; F = "Helse Sør-Øst"; G = 111919 },
    @{ Row = 25; A = 38; B = "Vestfold og Telemark";  C = 38; D = 8;  E = 1; F = "Helse Sør-Øst"; G = 111919 },
    @{ Row = 26; A = 42; B = "Agder";                 C = 42; D = 7;  E = 1; F = "Helse Sør-Øst"; G = 111919 },
    @{ Row = 27; A = 46; B = "Vestland";              C = 46; D = 5;  E = 2; F = "Helse Vest";    G = 100021 },
    @{ Row = 28; A = 54; B = "Troms og Finnmark";     C = 54; D = 1;  E = 4; F = "Helse Nord";     G = 100022 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

# --- Restore the selected cell that Excel recorded on save.
$ws.Range("P6").Select()
